# edit.ps1 - Apply updated cryptocurrency price/volume data to cryptos.xlsx
# Commit: "Updated cryptos list on Thu Aug 10 06:19:43 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell's value to be stored as TEXT (not auto-coerced to a
# number) by temporarily switching the cell to a text number format, the
# same trick a human editor uses in the Excel UI ( Format Cells > Text )
# before typing a value like "0.9989" or "1.000" that would otherwise be
# parsed as a numeric literal.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Rows whose Coin/Link didn't change: update Price (D) / Volume(1h) (E) only ---
$ws.Range("D2").Value = "29.545.69"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.852.59"
$ws.Range("E3").Value = "  -0.15%  "
Set-TextValue "D4" "0.9989"
Set-TextValue "D5" "243.13"
$ws.Range("E5").Value = "  -0.60%  "
Set-TextValue "D6" "0.6336"
$ws.Range("E6").Value = "  -1.25%  "
Set-TextValue "D7" "0.9999"
$ws.Range("E7").Value = "  -0.03%  "
Set-TextValue "D8" "47.86"
$ws.Range("E8").Value = "  +1.35%  "
Set-TextValue "D9" "0.07569"
$ws.Range("E9").Value = "  +1.14%  "
Set-TextValue "D10" "0.2981"
$ws.Range("E10").Value = "  +0.48%  "
Set-TextValue "D11" "24.33"
$ws.Range("E11").Value = "  -0.27%  "
Set-TextValue "D12" "0.07685"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.881.85"
$ws.Range("E13").Value = "  +1.35%  "
Set-TextValue "D14" "5.036"
$ws.Range("E14").Value = "  +0.00%  "
Set-TextValue "D15" "0.6879"
$ws.Range("E15").Value = "  -0.31%  "
Set-TextValue "D16" "83.95"
$ws.Range("E16").Value = "  +0.04%  "
Set-TextValue "D17" "0.000009864"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "2.132.82"
$ws.Range("E18").Value = "  +1.03%  "
Set-TextValue "D19" "6.241"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").Value = "29.592.39"
$ws.Range("E20").Value = "  -0.50%  "
Set-TextValue "D21" "235.35"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  -0.03%  "
Set-TextValue "D24" "7.650"
$ws.Range("E24").Value = "  +2.48%  "
Set-TextValue "D25" "1.000"
$ws.Range("E25").Value = "  -0.02%  "
Set-TextValue "D26" "155.95"
$ws.Range("E26").Value = "  -1.52%  "
Set-TextValue "D27" "0.1389"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -0.82%  "
Set-TextValue "D30" "1.484"
$ws.Range("E30").Value = "  -0.94%  "
Set-TextValue "D31" "0.05842"
$ws.Range("E31").Value = "  -5.74%  "
Set-TextValue "D32" "1.283"
$ws.Range("E32").Value = "  +0.94%  "
Set-TextValue "D33" "4.119"
$ws.Range("E33").Value = "  -0.69%  "
Set-TextValue "D34" "4.054"
$ws.Range("E34").Value = "  -0.87%  "
Set-TextValue "D35" "1.899"
$ws.Range("E35").Value = "  +0.27%  "
Set-TextValue "D36" "1.172"
$ws.Range("E36").Value = "  +0.05%  "
Set-TextValue "D37" "0.7195"
$ws.Range("E37").Value = "  -1.05%  "
Set-TextValue "D38" "2.594"
$ws.Range("E38").Value = "  -0.80%  "
Set-TextValue "D41" "0.01776"
$ws.Range("E41").Value = "  -0.39%  "
Set-TextValue "D42" "0.9154"
$ws.Range("E42").Value = "  -0.71%  "
Set-TextValue "D43" "6.144"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "2.040.21"
$ws.Range("E44").Value = "  +1.09%  "
Set-TextValue "D45" "0.9994"
$ws.Range("E45").Value = "  -0.12%  "
Set-TextValue "D48" "7.377"
$ws.Range("E48").Value = "  +10.42%  "
Set-TextValue "D49" "9.187"
$ws.Range("E49").Value = "  +0.09%  "
Set-TextValue "D50" "0.4047"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  -2.81%  "

# --- Rows that swapped rank order: Coin, Link, Price and Volume all change ---
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D39" "2.803"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.238.06"
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "67.71"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D47" "101.99"
$ws.Range("E47").Value = "  -0.24%  "
